# fix(publipostage): Refactor synthetic array /3
# Replace the colored-square emoji icons in column A with book emoji,
# and rename the "noir" (black) label to "bleu" (blue) to match the
# new blue-book icon.
#
# Note: cell values are written directly (by known address) rather than
# read-back-and-compared, since round-tripping non-BMP emoji through
# Range.Value reads is unreliable in this runtime; direct assignment of
# emoji string literals works correctly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: 🟥 -> 📕, ⬛ -> 📘, 🟧 -> 📙, 🟩 -> 📗
$ws.Range("A2").Value = "📕"
$ws.Range("A3").Value = "📕"
$ws.Range("A4").Value = "📘"
$ws.Range("A5").Value = "📕"
$ws.Range("A6").Value = "📕"
$ws.Range("A7").Value = "📙"
$ws.Range("A8").Value = "📗"
$ws.Range("A9").Value = "📗"

# Column B: "noir" -> "bleu" (row 4, matching the ⬛/📘 icon)
$ws.Range("B4").Value = "bleu"
